$d = $word.ActiveDocument

$replacements = @(
    @("981÷5=", "508÷3="),
    @("340÷2=", "789÷6="),
    @("289÷4=", "762÷3="),
    @("339÷2=", "714÷2="),
    @("942÷6=", "612÷2="),
    @("503÷5=", "594÷3="),
    @("972÷2=", "526÷2="),
    @("672÷9=", "699÷5="),
    @("862÷7=", "818÷6="),
    @("793÷9=", "641÷9="),
    @("717÷6=", "317÷2="),
    @("113÷2=", "169÷7="),
    @("633÷3=", "816÷9="),
    @("936÷4=", "605÷6="),
    @("418÷9=", "665÷9="),
    @("768÷9=", "408÷2="),
    @("391÷3=", "915÷4="),
    @("277÷5=", "146÷3="),
    @("848÷3=", "563÷4="),
    @("889÷5=", "133÷2="),
    @("574÷8=", "659÷8="),
    @("668÷9=", "751÷3="),
    @("552÷6=", "738÷6="),
    @("816÷4=", "791÷5="),
    @("229÷4=", "146÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
